$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 5 new "BD" (birthday) events for HelloLL members, continuing the
# existing event log (rows 2-98) with rows 99-103.

# Row 99: Hazuki Ren
$ws.Range("A98").Copy()
$ws.Range("A99").PasteSpecial(-4122)
$ws.Range("A99").Value = 40506
$ws.Range("B99").Value = "BD"
$ws.Range("C99").Value = "Hazuki Ren"

# Row 100: Heanna Sumire
$ws.Range("A98").Copy()
$ws.Range("A100").PasteSpecial(-4122)
$ws.Range("A100").Value = 40449
$ws.Range("B100").Value = "BD"
$ws.Range("C100").Value = "Heanna Sumire"

# Row 101: Arashi Chisato
$ws.Range("A98").Copy()
$ws.Range("A101").PasteSpecial(-4122)
$ws.Range("A101").Value = 40234
$ws.Range("B101").Value = "BD"
$ws.Range("C101").Value = "Arashi Chisato"

# Row 102: Thảng Khửa Khừa
$ws.Range("A98").Copy()
$ws.Range("A102").PasteSpecial(-4122)
$ws.Range("A102").Value = 40366
$ws.Range("B102").Value = "BD"
$ws.Range("C102").Value = "Thảng Khửa Khừa"

# Row 103: Shibuya Kanon
$ws.Range("A98").Copy()
$ws.Range("A103").PasteSpecial(-4122)
$ws.Range("A103").Value = 40299
$ws.Range("B103").Value = "BD"
$ws.Range("C103").Value = "Shibuya Kanon"

# Match the saved view state: scrolled down with C104 selected.
$ws.Range("C104").Select()
